$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 / B38 -- rewritten introduction summary
$ws.Cells.Item(38, 2).Value2 = "We have rewritten the introduction to more clearly expose the`nmotivation and to better captivate the audience. The words “this project” no longer appear. The introduction now begins with a general commentary on how different environmental effects complicate analysis and how the different effects leaves designers decreasing impacts in one dimension while increasing it another one. It then proceeds to explain how designers can use weighting techniques to assess multiple impacts at once and how this technique complicates assessments. Finally, the introduction presents a specific pair of products to use in a case study on how weighting might dictate design practices."
$ws.Rows.Item(38).RowHeight = 141.7

# Row 40 / B40 -- weighting triangle justification
$ws.Cells.Item(40, 2).Value2 = "The introduction now contains a more elaborate justification for the`nweighting triangle on page 2 where it explains how the study concerns the most common analysis practices and how SimaPro, which utilizes the weighting triangle, is one of the two most common software packages used. We do not use tradition to try to justify anything."

# Row 42 / B42 -- LCIA clarification (rich text, 'not' stays italic)
$ws.Cells.Item(42, 2).Value2 = "We have reduced the number of introductions of the weighting factor.`nAlso, the introduction now contains a clarification of the LCIA on page 2 We stress that LCIA does not amount to  LCA + Eco-indicator 99. It only refers to an axuliary phase of the LCA, and Eco-Indicator 99 is one  method available for implementing the LCIA. We believe the new explanation makes clearer what LCIA is and how it relates to the overall LCA and Eco-indicator 99."
$c42a = $ws.Cells.Item(42, 2).Characters(1, 168)
$c42a.Font.Name = "Arial"
$c42a.Font.Size = 12
$c42a.Font.Color = 0
$c42b = $ws.Cells.Item(42, 2).Characters(169, 3)
$c42b.Font.Name = "Arial"
$c42b.Font.Size = 12
$c42b.Font.Color = 0
$c42b.Font.Italic = $true
$c42c = $ws.Cells.Item(42, 2).Characters(172, 274)
$c42c.Font.Name = "Arial"
$c42c.Font.Size = 12
$c42c.Font.Color = 0

# Row 83 / B83 -- LCA and Eco-Indicator 99 justification
$ws.Cells.Item(83, 2).Value2 = "The introduction now has a justification for LCA and Eco-Indicator 99`nbased on their prevalence and our intention to focus on common practices. We emphasize the data limitations to make the reader aware of the study’s limitations, but the methodologies in this paper need not use these particular data. We use these data because of their availability given our limited resources, but they are not inherent to our methodogies. However, we concede that the data introduce uncertainties to our results. We could not eliminate all of the limitations in our work."
$ws.Rows.Item(83).RowHeight = 116.2

# Row 85 / B85 -- abstract no longer mentions cfls
$ws.Cells.Item(85, 2).Value2 = "The abstract no longer mentions cfls, which now focuses on the paper’s central case study. We only mention the cfls in the introduction as an example of the tradeoffs in green technology."
$ws.Rows.Item(85).RowHeight = 39.55

# Row 87 / B87 -- use hours reference
$ws.Cells.Item(87, 2).Value2 = "The use hours reference now appears in the same sentence in which`nwe first reference use hours, “The desktop computers have more use hours because, in practice, users usually switch off thin clients at night, but only 30% of users switch off desktop computers [2].”"
$ws.Rows.Item(87).RowHeight = 77.95

# Row 124 / D124 -- our own numbers / no realistic way
$ws.Cells.Item(124, 4).Value2 = "We use information from databases.`nI think the reviewer has suggested`nthat we get our own numbers from`nmeasuring material instead of using`nthe databases. I see no realistic`nway to satisfy this recommendation."

# Row 126 / D126 -- documentation wording
$ws.Cells.Item(126, 4).Value2 = "I suppose we could add a section`nthat shows how other researchers`nhave conducted LCAs, but I don’t`nknow if it would add anything to the`npaper. I didn’t do anything unique`nwith the LCA and followed the`nofficial documentation as closely as I could, so I don’t think we have`nanything to say about how our`nmethod differs from existing LCAs."

# Sheet view: scroll position & selection
$excel.ActiveWindow.ScrollRow = 121
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B126").Select()
